# "working on Storage page"
#
# - Move the active/selected tab from "ValidCredAndCaptcha" to
#   "ValidCredAndInValidCaptcha" (last sheet): this both clears
#   tabSelected on the former and sets it on the latter, and updates the
#   workbook's activeTab.
# - Update ValidCredAndInValidCaptcha's selection to E2 (single cell)
#   instead of D2:E2.
# - Update the shared "Password" link text from "ipmcloud@300#" to
#   "ipmcloud@1#" (used by the Valid / ValidCredAndCaptcha /
#   ValidCredAndInValidCaptcha sheets).

$wb = $excel.ActiveWorkbook

$validSheet = $wb.Worksheets.Item("Valid")
$validCredAndCaptchaSheet = $wb.Worksheets.Item("ValidCredAndCaptcha")
$validCredAndInValidCaptchaSheet = $wb.Worksheets.Item("ValidCredAndInValidCaptcha")

# Update the password text everywhere it is used (shared string, so one
# write per sheet keeps every usage in sync).
$validSheet.Range("C2").Value = "ipmcloud@1#"
$validCredAndCaptchaSheet.Range("E2").Value = "ipmcloud@1#"
$validCredAndInValidCaptchaSheet.Range("E2").Value = "ipmcloud@1#"

# Move the selection on the last sheet from D2:E2 to just E2.
$validCredAndInValidCaptchaSheet.Range("E2").Select()

# Activate the last sheet so it becomes the selected/active tab
# (tabSelected moves off ValidCredAndCaptcha and onto this sheet, and the
# workbook's activeTab updates accordingly).
$validCredAndInValidCaptchaSheet.Activate()
